# "make it so that each person can only do 1 thing per month"
# Remove extra "X" marks so each volunteer is assigned at most one
# role/task, update the "Max weeks per month" counts to match, and
# move the active-cell selection to where the last edit landed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Kyle (row 2): keep only DJ (E2); drop Opening/Teaching lead/Teaching follow ---
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""

# --- Colby (row 4): drop Opening and First Door Shift ---
$ws.Range("B4").Value = ""

# --- Geoff (row 7): Max weeks per month 2 -> 3 ---
$ws.Range("J7").Value = 3

# --- Madeline (row 10): Max weeks per month 1 -> 2, row height shrinks ---
$ws.Range("J10").Value = 2
$ws.Rows.Item(10).RowHeight = 15.75

# --- Jessica (row 11): add Closing ---
$ws.Range("F11").Value = "X"

# --- Alex (row 12): drop Opening/Teaching lead/Teaching follow/DJ/Closing, keep moose member + Promotion ---
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("J12").Value = 0

# --- move selection to reflect where editing finished ---
[void]$ws.Range("J11").Select()
